$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 775 (shifts existing rows 775:816 down to 776:817)
$ws.Rows(775).Insert()

# Populate the newly inserted row with the new record.
# Force the date column to stay plain text (matching the rest of column A)
# instead of being auto-parsed into a date serial number by COM.
$ws.Range("A775").NumberFormat = "@"
$ws.Range("A775").Value = "2026/02/05"
$ws.Range("A775").ClearFormats()

$ws.Range("B775").Value = "木"
$ws.Range("C775").Value = 14
$ws.Range("D775").Value = 201
